# Apply updated cryptocurrency price/volume data and reorder two coin pairs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.704.81'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.919.04'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.96'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4938'
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3010'
$ws.Range("E8").Value = '  +2.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06783'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = '1.908.01'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.24'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07332'
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.234'
$ws.Range("E13").Value = '  +3.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.68'
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '30.685.88'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007982'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.57'
$ws.Range("E18").Value = '  +3.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '2.168.82'
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.385'
$ws.Range("E21").Value = '  +11.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '197.21'
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.333'
$ws.Range("E24").Value = '  +4.31%  '
$ws.Range("E25").Value = '  +3.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.08'
$ws.Range("E26").Value = '  +4.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.64'
$ws.Range("E27").Value = '  -2.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.964'
$ws.Range("E28").Value = '  +3.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.472'
$ws.Range("E29").Value = '  +5.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.364'
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09151'
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.088'
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05276'
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7449'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.709'
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("E37").Value = '  +1.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.726'
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9278'
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.089'
$ws.Range("E40").Value = '  -2.61%  '
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.98'
$ws.Range("E42").Value = '  +24.42%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.11'
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.939'
$ws.Range("E44").Value = '  +3.47%  '
$ws.Range("E45").Value = '  +4.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.730'
$ws.Range("E47").Value = '  +1.74%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.069'
$ws.Range("E48").Value = '  +4.10%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.38'
$ws.Range("E49").Value = '  +5.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05892'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("E51").Value = '  +3.17%  '
